$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test CU 1 Devenir membre")

# Row 8: change result from "Ko" (red) to "Ok" (green), matching style of B2:B7,
# and clear the "Fonctionnalité non implémentée" observation in C8.
$ws.Range("B8").Value = "Ok"
$ws.Range("B8").Font.Color = $ws.Range("B2").Font.Color

$ws.Range("C8").ClearContents()

# Move the active selection to C8
$ws.Activate()
$ws.Range("C8").Select()
